# Apply updated dSF (column F) values as described by the commit:
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3"  = -5
    "F6"  = -2
    "F7"  = 7
    "F12" = -6
    "F16" = 0
    "F17" = 2
    "F19" = -2
    "F20" = -4
    "F22" = -1
    "F25" = -2
    "F29" = -2
    "F30" = -2
    "F31" = -4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
